# Update countries & provincias Spain
# Refresh the COVID figures pull and re-rank a handful of countries whose
# updated totals change their position in the (descending, by total cases)
# table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last refreshed" timestamp shown in A1 ---------------------
$ws.Range("A1").Value = "Datos actualizados a 23 de Julio de 2020 a las 06:57"

# --- Helper to rewrite one data row (columns A..H) --------------------------
function Set-Row($row, $country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Row 15 - Pakistan: refreshed totals (no reorder)
Set-Row 15 "Pakistan" 269191 1763 213175 50307 0 32 5709

# Row 36 - Belgica: refreshed totals (no reorder)
Set-Row 36 "Belgica" 64627 369 17347 37472 0 3 9808

# Rows 56/57 - Kirguistan overtakes Ghana
Set-Row 56 "Kirguistan" 30349 990 16791 12389 0 46 1169
Set-Row 57 "Ghana" 29672 0 26090 3429 0 0 153

# Rows 106/107/108 - Malaui overtakes Tailandia and Somalia
Set-Row 106 "Malaui" 3302 153 1282 1944 0 5 76
Set-Row 107 "Tailandia" 3269 8 3105 106 0 0 58
Set-Row 108 "Somalia" 3161 0 1495 1573 0 0 93

# Rows 131/132/133 - Benin overtakes Ruanda and Yemen
Set-Row 131 "Benin" 1690 0 782 874 0 0 34
Set-Row 132 "Ruanda" 1689 0 867 817 0 0 5
Set-Row 133 "Yemen" 1640 0 751 431 0 0 458

# Row 170 - Mongolia: refreshed totals (no reorder)
$ws.Cells.Item(170, 4).Value = 214
$ws.Cells.Item(170, 5).Value = 73

# Row 175 - Camboya: refreshed totals (no reorder)
$ws.Cells.Item(175, 2).Value = 198
$ws.Cells.Item(175, 3).Value = 1
$ws.Cells.Item(175, 4).Value = 142
$ws.Cells.Item(175, 5).Value = 56
